$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "duplicate_image_filename" column (E) was missing values for the
# practice/generic/unique_video/unique_audio stimulus rows (rows 2-21).
# Fill those in with "NA".
$ws.Range("E2:E21").Value = "NA"
